$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (2025-07-17 / BEMOL S/A / 384275 / KIT LANCHE...)
# this shifts all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# A handful of cells (beyond the pure shift) carry updated values in the
# refreshed data pull - fix those up on the now-shifted rows.
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = -19
$ws.Range("H3").Value = 1.07
$ws.Range("I3").Value = 0.27
$ws.Range("G4").Value = -247
$ws.Range("G5").Value = -108
